# Add a new "mCoffee" folder entry to the file-listing table on Sheet1
# and re-sort the table (Type, then Name) so the new row lands in its
# alphabetically-correct position, exactly like the author did after
# adding a row and re-applying the sheet's existing sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently occupies A1:C43 (header in row 1). Insert a fresh
# row right under the header (any position works - we sort afterwards)
# and populate it with the new file/folder description.
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(6, 1).Value2 = "mCoffee"
$ws.Cells.Item(6, 2).Value2 = "folder"
$ws.Cells.Item(6, 3).Value2 = "mCoffee alignments of the nonDBD (non) and full-length proteins (fl) for all paralog pairs and their orthologs "

# Re-sort the whole table (now A1:C44) by Type then Name, matching the
# workbook's existing sortState (B then A), so the new row is placed in
# its correct alphabetical spot along with everything else.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()

$key1 = $ws.Range("B2:B44")
$key2 = $ws.Range("A2:A44")
$sortObj.SortFields.Add($key1)
$sortObj.SortFields.Add($key2)

$sortObj.SetRange($ws.Range("A1:C44"))
$sortObj.Header = 1
$sortObj.Apply()

# Reflect the new extent in the sheet's selection (the table grew by one
# row, from C44 down to the whole A1:C44 block).
$ws.Range("A1:C44").Select()
